# Data update using git
# Updates the "Pagos" (F) and "Inscrições homologadas" (H) columns on the
# "Inscricoes" sheet for the rows whose payment counts increased.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("F3").Value = 45
$ws.Range("H3").Value = 46

$ws.Range("F5").Value = 141
$ws.Range("H5").Value = 152

$ws.Range("F6").Value = 51
$ws.Range("H6").Value = 61

$ws.Range("F10").Value = 634
$ws.Range("H10").Value = 729

$ws.Range("F11").Value = 413
$ws.Range("H11").Value = 478

$ws.Range("F12").Value = 671
$ws.Range("H12").Value = 757

$ws.Range("F13").Value = 141
$ws.Range("H13").Value = 175

$ws.Range("F16").Value = 180
$ws.Range("H16").Value = 228

$ws.Range("F23").Value = 174
$ws.Range("H23").Value = 226

$ws.Range("F24").Value = 246
$ws.Range("H24").Value = 276

$ws.Range("F25").Value = 285
$ws.Range("H25").Value = 345

$ws.Range("F26").Value = 211
$ws.Range("H26").Value = 236

$ws.Range("F27").Value = 314
$ws.Range("H27").Value = 396

$ws.Range("F28").Value = 176
$ws.Range("H28").Value = 228

$ws.Range("F30").Value = 222
$ws.Range("H30").Value = 275

$ws.Range("F38").Value = 93
$ws.Range("H38").Value = 110

$ws.Range("F41").Value = 353
$ws.Range("H41").Value = 445

$ws.Range("F42").Value = 415
$ws.Range("H42").Value = 476

$ws.Range("F44").Value = 315
$ws.Range("H44").Value = 383

$ws.Range("F45").Value = 158
$ws.Range("H45").Value = 197

$ws.Range("F47").Value = 474
$ws.Range("H47").Value = 566

$ws.Range("F48").Value = 209
$ws.Range("H48").Value = 253
